# Slide Deck 1 (e2immu) - April 2021 content refresh
# 1) Title slide: "E2IMMU" -> "e2immu" (italic)
# 2) "What is E2IMMU?" -> "What is e2immu?" (e2immu italic)
# 3) "Does E2IMMU do typical analyser stuff?" -> "Does e2immu do typical analyser stuff?" (e2immu italic)
# 4) "The majority of classes..." -> "The majority of types..."

$p = $ppt.ActivePresentation

# --- Slide 1 (position 1): Title "Title 1" shape, id=2 ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$run1 = $tr1.Characters(1, 6)
$run1.Text = "e2immu"
$run1.Font.Italic = $true

# --- Slide 2 (position 2): Title "Title 1" shape, id=2: "What is E2IMMU?" ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(1)
$tr2 = $sh2.TextFrame.TextRange
$run2 = $tr2.Characters(9, 6)
$run2.Text = "e2immu"
$run2.Font.Italic = $true

# --- Slide 6 (position 6): Title "Title 1" shape, id=2: "Does E2IMMU do typical analyser stuff?" ---
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(1)
$tr6 = $sh6.TextFrame.TextRange
$run6 = $tr6.Characters(6, 6)
$run6.Text = "e2immu"
$run6.Font.Italic = $true

# --- Slide 10 (position 10): Content Placeholder 2, id=3: last paragraph "classes" -> "types" ---
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$fullText10 = $tr10.Text
$lastBreak = $fullText10.LastIndexOf([char]13)
$paraStart = $lastBreak + 2
$paraLen = $fullText10.Length - $lastBreak - 1
$lastPara = $tr10.Characters($paraStart, $paraLen)
$lastPara.Text = "The majority of types in your project should be containers"

Write-Output "done"
